# Login Page Test script with data provider
#
# - Rename the second worksheet ("Sheet1") to "LoginData" and populate it
#   with a small login data-provider table (UserEmail/Password header +
#   one row of sample credentials).
# - Make "LoginData" the active sheet/tab (was "Data" before).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(2)
$ws.Name = "LoginData"

$ws.Range("A1").Value = "UserEmail"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "automation101@gmail.com"
$ws.Range("B2").Value = "automation101"

$ws.Columns.Item(1).ColumnWidth = 11.42

# Activating the sheet flips the workbook's activeTab / each sheet's
# tabSelected flag, and records the new selection/activeCell.
$ws.Activate()
$ws.Range("H19").Select() | Out-Null
